$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update country names (shared string text) that were reordered/relabeled
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 9 de Abril de 2020 a las 20:22'
$ws.Cells.Item(6, 1).Value = 'Italia'
$ws.Cells.Item(7, 1).Value = 'Francia'
$ws.Cells.Item(8, 1).Value = 'Alemania'
$ws.Cells.Item(9, 1).Value = 'China'
$ws.Cells.Item(161, 1).Value = 'Tanzania'
$ws.Cells.Item(162, 1).Value = 'Libia'
$ws.Cells.Item(163, 1).Value = 'Birmania'
$ws.Cells.Item(164, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(185, 1).Value = 'Seychelles'
$ws.Cells.Item(186, 1).Value = 'Republica del Chad'
$ws.Cells.Item(187, 1).Value = 'Zimbabue'
$ws.Cells.Item(188, 1).Value = 'Groenlandia'
$ws.Cells.Item(189, 1).Value = 'Surinam'
$ws.Cells.Item(194, 1).Value = 'San Vicente y las Granadinas'
$ws.Cells.Item(195, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(196, 1).Value = 'Malaui'
$ws.Cells.Item(197, 1).Value = 'Santa Sede'
$ws.Cells.Item(198, 1).Value = 'Sierra Leona'
$ws.Cells.Item(199, 1).Value = 'Nicaragua'
$ws.Cells.Item(201, 1).Value = 'Cabo Verde'
$ws.Cells.Item(202, 1).Value = 'San Bartolome'
$ws.Cells.Item(203, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(209, 1).Value = 'Burundi'
$ws.Cells.Item(210, 1).Value = 'Anguila'
$ws.Cells.Item(211, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(212, 1).Value = 'Bonaire, San Eustaquio y Saba'
$ws.Cells.Item(215, 1).Value = 'San Pedro y Miquelon'

# Update numeric statistics cells
$ws.Cells.Item(4, 2).Value = 455445
$ws.Cells.Item(4, 3).Value = 20518
$ws.Cells.Item(4, 5).Value = 414769
$ws.Cells.Item(4, 6).Value = 9813
$ws.Cells.Item(4, 7).Value = 1326
$ws.Cells.Item(4, 8).Value = 16114
$ws.Cells.Item(7, 2).Value = 117749
$ws.Cells.Item(7, 3).Value = 4799
$ws.Cells.Item(7, 4).Value = 23206
$ws.Cells.Item(7, 5).Value = 82333
$ws.Cells.Item(7, 6).Value = 7066
$ws.Cells.Item(7, 7).Value = 1341
$ws.Cells.Item(7, 8).Value = 12210
$ws.Cells.Item(8, 2).Value = 115523
$ws.Cells.Item(8, 3).Value = 2227
$ws.Cells.Item(8, 4).Value = 46300
$ws.Cells.Item(8, 5).Value = 66772
$ws.Cells.Item(8, 6).Value = 4895
$ws.Cells.Item(8, 7).Value = 102
$ws.Cells.Item(8, 8).Value = 2451
$ws.Cells.Item(19, 2).Value = 13237
$ws.Cells.Item(19, 3).Value = 295
$ws.Cells.Item(19, 5).Value = 7702
$ws.Cells.Item(24, 2).Value = 6725
$ws.Cells.Item(24, 3).Value = 809
$ws.Cells.Item(24, 5).Value = 5924
$ws.Cells.Item(42, 4).Value = 374
$ws.Cells.Item(42, 5).Value = 2689
$ws.Cells.Item(90, 5).Value = 345
$ws.Cells.Item(90, 7).Value = 1
$ws.Cells.Item(90, 8).Value = 3
$ws.Cells.Item(92, 2).Value = 443
$ws.Cells.Item(92, 3).Value = 29
$ws.Cells.Item(92, 4).Value = 146
$ws.Cells.Item(92, 5).Value = 273
$ws.Cells.Item(92, 7).Value = 1
$ws.Cells.Item(92, 8).Value = 24
$ws.Cells.Item(101, 4).Value = 33
$ws.Cells.Item(101, 5).Value = 276
$ws.Cells.Item(150, 4).Value = 24
$ws.Cells.Item(150, 5).Value = 14
$ws.Cells.Item(162, 2).Value = 24
$ws.Cells.Item(162, 3).Value = 3
$ws.Cells.Item(162, 4).Value = 8
$ws.Cells.Item(162, 5).Value = 15
$ws.Cells.Item(162, 8).Value = 1
$ws.Cells.Item(163, 2).Value = 23
$ws.Cells.Item(163, 3).Value = 1
$ws.Cells.Item(163, 4).Value = 1
$ws.Cells.Item(163, 5).Value = 19
$ws.Cells.Item(163, 8).Value = 3
$ws.Cells.Item(186, 3).Value = 1
$ws.Cells.Item(186, 4).Value = 2
$ws.Cells.Item(186, 5).Value = 9
$ws.Cells.Item(186, 8).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 8
$ws.Cells.Item(187, 8).Value = 3
$ws.Cells.Item(188, 2).Value = 11
$ws.Cells.Item(188, 4).Value = 11
$ws.Cells.Item(188, 5).Value = 0
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(196, 6).Value = 1
$ws.Cells.Item(199, 3).Value = 1
$ws.Cells.Item(199, 4).Value = 0
$ws.Cells.Item(199, 5).Value = 6
$ws.Cells.Item(200, 3).Value = 0
$ws.Cells.Item(200, 4).Value = 1
$ws.Cells.Item(200, 5).Value = 5
$ws.Cells.Item(201, 2).Value = 7
$ws.Cells.Item(201, 3).Value = 1
$ws.Cells.Item(201, 4).Value = 2
$ws.Cells.Item(201, 5).Value = 4
$ws.Cells.Item(201, 8).Value = 1
$ws.Cells.Item(202, 4).Value = 1
$ws.Cells.Item(202, 8).Value = 0
